# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder "Periodo Mora" column (E16:E20) so periods run ascending 2402..2406
$ws.Range("E16").Value = "2402"
$ws.Range("E17").Value = "2403"
$ws.Range("E18").Value = "2404"
$ws.Range("E19").Value = "2405"
$ws.Range("E20").Value = "2406"

# Update "Valor Mora" (F16:F20) to match the new period ordering
$ws.Range("F16").Value = 46400
$ws.Range("F17").Value = 46400
$ws.Range("F18").Value = 46400
$ws.Range("F19").Value = 52000
$ws.Range("F20").Value = 20800

# Update "Salario Basico" (G16:G20) for every row
$ws.Range("G16:G20").Value = 1160000
